$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2426.5881
$ws.Range("I64").Value = 2246.5715
$ws.Range("J64").Value = 3266.6667
$ws.Range("K64").Value = 2246.5715
$ws.Range("L64").Value = 3266.6667
$ws.Range("M64").Value = -1998.5715
$ws.Range("N64").Value = -3762.6667
$ws.Range("H67").Value = 2426.5881
$ws.Range("I67").Value = 2246.5715
$ws.Range("J67").Value = 3266.6667
$ws.Range("K67").Value = 2246.5715
$ws.Range("L67").Value = 3266.6667
$ws.Range("M67").Value = -1388.5715
$ws.Range("N67").Value = -4982.6667
$ws.Range("H132").Value = 3835.7964
$ws.Range("I132").Value = 3693.9556
$ws.Range("K132").Value = 11081.8668
$ws.Range("M132").Value = -8551.8668
$ws.Range("H135").Value = 598
$ws.Range("I135").Value = 397.7143
$ws.Range("K135").Value = 3579.4287
$ws.Range("M135").Value = -1044.4287
$ws.Range("H137").Value = 2115.898
$ws.Range("I137").Value = 1074.1351
$ws.Range("J137").Value = 5328
$ws.Range("K137").Value = 3222.4053
$ws.Range("L137").Value = 15984
$ws.Range("M137").Value = -672.4052999999999
$ws.Range("N137").Value = -21084
$ws.Range("H138").Value = 3687.59
$ws.Range("I138").Value = 1508.8462
$ws.Range("J138").Value = 4013.1494
$ws.Range("K138").Value = 4526.5386
$ws.Range("L138").Value = 12039.4482
$ws.Range("M138").Value = 613.4614000000001
$ws.Range("N138").Value = -22319.4482
$ws.Range("H141").Value = 8160.1
$ws.Range("I141").Value = 9434.708000000001
$ws.Range("J141").Value = 3061.6667
$ws.Range("K141").Value = 28304.124
$ws.Range("L141").Value = 9185.000100000001
$ws.Range("M141").Value = -23124.124
$ws.Range("N141").Value = -19545.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4017.6052
$ws.Range("I32").Value = 3105.7715
$ws.Range("J32").Value = 14655.667
$ws.Range("K32").Value = 3105.7715
$ws.Range("L32").Value = 14655.667
$ws.Range("M32").Value = -2818.7715
$ws.Range("N32").Value = -15229.667
$ws.Range("H61").Value = 982.4286
$ws.Range("I61").Value = 802.8
$ws.Range("K61").Value = 802.8
$ws.Range("M61").Value = -590.8
$ws.Range("H64").Value = 48378.832
$ws.Range("J64").Value = 48378.832
$ws.Range("L64").Value = 48378.832
$ws.Range("N64").Value = -48874.832
$ws.Range("H67").Value = 48378.832
$ws.Range("J67").Value = 48378.832
$ws.Range("L67").Value = 48378.832
$ws.Range("N67").Value = -50094.832
$ws.Range("H68").Value = 50099
$ws.Range("J68").Value = 50099
$ws.Range("L68").Value = 50099
$ws.Range("N68").Value = -51721
$ws.Range("H71").Value = 50099
$ws.Range("J71").Value = 50099
$ws.Range("L71").Value = 150297
$ws.Range("N71").Value = -158409
$ws.Range("H74").Value = 3029.282
$ws.Range("I74").Value = 3180.1
$ws.Range("J74").Value = 2526.5557
$ws.Range("K74").Value = 3180.1
$ws.Range("L74").Value = 2526.5557
$ws.Range("M74").Value = -2306.1
$ws.Range("N74").Value = -4274.5557
$ws.Range("H77").Value = 3029.282
$ws.Range("I77").Value = 3180.1
$ws.Range("J77").Value = 2526.5557
$ws.Range("K77").Value = 15900.5
$ws.Range("L77").Value = 12632.7785
$ws.Range("M77").Value = -11532.5
$ws.Range("N77").Value = -21368.7785
$ws.Range("H80").Value = 25322
$ws.Range("J80").Value = 26586.4
$ws.Range("L80").Value = 26586.4
$ws.Range("N80").Value = -28582.4
$ws.Range("H83").Value = 25322
$ws.Range("J83").Value = 26586.4
$ws.Range("L83").Value = 79759.20000000001
$ws.Range("N83").Value = -89743.20000000001
$ws.Range("H110").Value = 1610.9474
$ws.Range("I110").Value = 1607.2
$ws.Range("J110").Value = 1625
$ws.Range("K110").Value = 1607.2
$ws.Range("L110").Value = 1625
$ws.Range("M110").Value = 437.8
$ws.Range("N110").Value = -5715
$ws.Range("H122").Value = 2569.4285
$ws.Range("I122").Value = 1397.2
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 4191.6
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -1741.6
$ws.Range("N122").Value = -21400
$ws.Range("H136").Value = 982.4286
$ws.Range("I136").Value = 802.8
$ws.Range("K136").Value = 2408.4
$ws.Range("M136").Value = 141.6000000000004

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1284.8334
$ws.Range("I86").Value = 1256.5416
$ws.Range("J86").Value = 1398
$ws.Range("K86").Value = 1256.5416
$ws.Range("L86").Value = 1398
$ws.Range("M86").Value = -133.5416
$ws.Range("N86").Value = -3644
$ws.Range("H89").Value = 1284.8334
$ws.Range("I89").Value = 1256.5416
$ws.Range("J89").Value = 1398
$ws.Range("K89").Value = 6282.708000000001
$ws.Range("L89").Value = 6990
$ws.Range("M89").Value = -666.7080000000005
$ws.Range("N89").Value = -18222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10872034
$ws.Range("I31").Value = 1367.2812
$ws.Range("J31").Value = 35719270
$ws.Range("K31").Value = 1367.2812
$ws.Range("L31").Value = 35719270
$ws.Range("M31").Value = -1072.2812
$ws.Range("N31").Value = -35719860
$ws.Range("H34").Value = 10872034
$ws.Range("I34").Value = 1367.2812
$ws.Range("J34").Value = 35719270
$ws.Range("K34").Value = 1367.2812
$ws.Range("L34").Value = 35719270
$ws.Range("M34").Value = -1165.2812
$ws.Range("N34").Value = -35719674
$ws.Range("H58").Value = 2000.8788
$ws.Range("I58").Value = 1710.629
$ws.Range("J58").Value = 6499.75
$ws.Range("K58").Value = 1710.629
$ws.Range("L58").Value = 6499.75
$ws.Range("M58").Value = -1507.629
$ws.Range("N58").Value = -6905.75
$ws.Range("H134").Value = 3626.3618
$ws.Range("I134").Value = 3905.3225
$ws.Range("K134").Value = 11715.9675
$ws.Range("M134").Value = -9180.967500000001
$ws.Range("H136").Value = 2000.8788
$ws.Range("I136").Value = 1710.629
$ws.Range("J136").Value = 6499.75
$ws.Range("K136").Value = 5131.887
$ws.Range("L136").Value = 19499.25
$ws.Range("M136").Value = -2581.887
$ws.Range("N136").Value = -24599.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 969.0952
$ws.Range("I113").Value = 745.1177
$ws.Range("J113").Value = 1921
$ws.Range("K113").Value = 2235.3531
$ws.Range("L113").Value = 5763
$ws.Range("M113").Value = -65.35310000000027
$ws.Range("N113").Value = -10103
$ws.Range("H131").Value = 8197539.5
$ws.Range("J131").Value = 861.43396
$ws.Range("L131").Value = 2584.30188
$ws.Range("N131").Value = -12664.30188

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2629.25
$ws.Range("I126").Value = 2662.221
$ws.Range("J126").Value = 2002.8
$ws.Range("K126").Value = 7986.663
$ws.Range("L126").Value = 6008.4
$ws.Range("M126").Value = -5516.663
$ws.Range("N126").Value = -10948.4
$ws.Range("H132").Value = 2290.756
$ws.Range("I132").Value = 1356.1613
$ws.Range("J132").Value = 5188
$ws.Range("K132").Value = 4068.4839
$ws.Range("L132").Value = 15564
$ws.Range("M132").Value = -1538.4839
$ws.Range("N132").Value = -20624

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4180.5557
$ws.Range("I122").Value = 2479.125
$ws.Range("J122").Value = 6655.364
$ws.Range("K122").Value = 7437.375
$ws.Range("L122").Value = 19966.092
$ws.Range("M122").Value = -4987.375
$ws.Range("N122").Value = -24866.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1734.5264
$ws.Range("I136").Value = 494.86047
$ws.Range("J136").Value = 5542.0713
$ws.Range("K136").Value = 1484.58141
$ws.Range("L136").Value = 16626.2139
$ws.Range("M136").Value = 1065.41859
$ws.Range("N136").Value = -21726.2139
